# Update the "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages data pull and the one generated at commit 456a3b4.
#
# Sheet "展览" (index 1):
#   F2: 1091 -> 1094   南宁·第五届小蜜蜂动漫嘉年华
#   F4: 1701 -> 1713   南宁·草莓动漫节
#   F6: 190  -> 191    南宁·布谷鸟动漫展4th
#
# Sheet "全部类型" (index 4) carries the same three events, but at
# different row offsets:
#   F2: 1091 -> 1094
#   F4: 1701 -> 1713
#   F7: 190  -> 191

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 1094
$wsExhibition.Range("F4").Value = 1713
$wsExhibition.Range("F6").Value = 191

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1094
$wsAll.Range("F4").Value = 1713
$wsAll.Range("F7").Value = 191
